# Update odds data in row 2 of the active worksheet to reflect the latest
# FlashScore odds for the match (Slavia Sofia vs Lok. Sofia, 24/10/2024).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value  = 1.75
$ws.Range("I2").Value  = 4.75
$ws.Range("J2").Value  = 2.4
$ws.Range("K2").Value  = 2.05
$ws.Range("L2").Value  = 5.5
$ws.Range("M2").Value  = 1.07
$ws.Range("N2").Value  = 9
$ws.Range("O2").Value  = 1.4
$ws.Range("P2").Value  = 2.75
$ws.Range("Q2").Value  = 2.25
$ws.Range("R2").Value  = 1.62
$ws.Range("S2").Value  = 1.5
$ws.Range("T2").Value  = 2.5
$ws.Range("U2").Value  = 2.1
$ws.Range("V2").Value  = 1.67
$ws.Range("W2").Value  = 5.5
$ws.Range("Z2").Value  = 13
$ws.Range("AE2").Value = 21
$ws.Range("AI2").Value = 23
$ws.Range("AJ2").Value = 17
$ws.Range("AM2").Value = 51
$ws.Range("AO2").Value = 9.5
$ws.Range("AT2").Value = 2.5
$ws.Range("AW2").Value = 6.5
$ws.Range("AX2").Value = 29
$ws.Range("AY2").Value = 41
$ws.Range("BA2").Value = 151
$ws.Range("BB2").Value = 351
